$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Status moved from "In Translation" to "Ready for handoff" (report generated for handoff)
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Refresh the handoff timestamps
$overview.Range("G2").Value = "2016-08-16 14:59:37"
$zhcn.Range("H2").Value = "2016-08-16 14:59:33"
$dede.Range("H2").Value = "2016-08-16 14:59:37"

# Widen the Status columns to fit the new, longer text
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333332
$dede.Columns.Item(3).ColumnWidth = 16.333333333333332
